$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update column F (想去人数) values
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 38
$wsExhibit.Range("F5").Value = 3418
$wsExhibit.Range("F7").Value = 414
$wsExhibit.Range("F10").Value = 38
$wsExhibit.Range("F11").Value = 1266
$wsExhibit.Range("F13").Value = 1570
$wsExhibit.Range("F14").Value = 112

# Sheet "全部类型" (sheet4): update column F (想去人数) values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 38
$wsAll.Range("F5").Value = 3418
$wsAll.Range("F7").Value = 414
$wsAll.Range("F11").Value = 38
$wsAll.Range("F14").Value = 1266
$wsAll.Range("F16").Value = 1570
$wsAll.Range("F17").Value = 113
